# Add 5 new "numbered system message" rows just above the trailing
# catch-all row (A=2000000), pushing that row from 278 down to 283,
# and select C289 / scroll so row 265 is at the top - matching the
# author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The former last data row (2000000 / "[1]") shifts from row 278 to 283.
# Write it first at its new location, then fill rows 278-282 with the
# five newly added messages, each following the existing table's
# pattern: col A = 1000000 + n, col B = 0, col C = message text.

$ws.Cells.Item(283, 1).Value = 2000000
$ws.Cells.Item(283, 2).Value = 0
$ws.Cells.Item(283, 3).Value = "[1]"

$ws.Cells.Item(278, 1).Value = 1000279
$ws.Cells.Item(278, 2).Value = 0
$ws.Cells.Item(278, 3).Value = "The operation code '[1]' for '[2]' is invalid!"

$ws.Cells.Item(279, 1).Value = 1000280
$ws.Cells.Item(279, 2).Value = 0
$ws.Cells.Item(279, 3).Value = "The INSERT ID for '[1]' '[2]' could not be added!"

$ws.Cells.Item(280, 1).Value = 1000281
$ws.Cells.Item(280, 2).Value = 0
$ws.Cells.Item(280, 3).Value = "The DELETE ID for '[1]' '[2]' could not be added!"

$ws.Cells.Item(281, 1).Value = 1000282
$ws.Cells.Item(281, 2).Value = 0
$ws.Cells.Item(281, 3).Value = "The INSERT ID '[1]' and DETAIL ID '[2]' could not be deleted or a condition prevents it from deletion!"

$ws.Cells.Item(282, 1).Value = 1000283
$ws.Cells.Item(282, 2).Value = 0
$ws.Cells.Item(282, 3).Value = "The INSERT ID '[1]' and DETAIL ID '[2]' could not be deleted because a lock condition cannot be applied on a base64 field!"

# Match the author's final selection / scroll position.
$ws.Range("C289").Select() | Out-Null
